# Add schema='SBtab' to the ObjTables document- and table-level metadata
# strings, reorder tableFormat/id ahead of name/date, and bump the embedded
# timestamp (as described in the commit "adding schema to document, table
# metdata").

$wb = $excel.ActiveWorkbook

$docDate   = "2020-03-09 23:58:57"
$tableDate = "2020-03-09 23:58:57"

# id -> (sheet tab name, date stamp for that table's metadata line)
$tables = @(
    @{ Id = "Compartment";            Date = $tableDate },
    @{ Id = "Compound";                Date = $tableDate },
    @{ Id = "Definition";              Date = $tableDate },
    @{ Id = "Enzyme";                  Date = $tableDate },
    @{ Id = "FbcObjective";            Date = $tableDate },
    @{ Id = "Gene";                    Date = $tableDate },
    @{ Id = "Layout";                  Date = $tableDate },
    @{ Id = "Measurement";             Date = $tableDate },
    @{ Id = "PbConfig";                Date = $tableDate },
    @{ Id = "Position";                Date = $tableDate },
    @{ Id = "Protein";                 Date = $tableDate },
    @{ Id = "Quantity";                Date = $tableDate },
    @{ Id = "QuantityInfo";            Date = $tableDate },
    @{ Id = "QuantityMatrix";          Date = $tableDate },
    @{ Id = "Reaction";                Date = $tableDate },
    @{ Id = "ReactionStoichiometry";   Date = $tableDate },
    @{ Id = "Regulator";               Date = $tableDate },
    @{ Id = "Relation";                Date = $tableDate },
    @{ Id = "Relationship";            Date = $tableDate },
    @{ Id = "SparseMatrix";            Date = $tableDate },
    @{ Id = "SparseMatrixColumn";      Date = $tableDate },
    @{ Id = "SparseMatrixOrdered";     Date = $tableDate },
    @{ Id = "SparseMatrixRow";         Date = $tableDate },
    @{ Id = "StoichiometricMatrix";    Date = $tableDate },
    @{ Id = "rxnconContingencyList";   Date = $tableDate },
    @{ Id = "rxnconReactionList";      Date = "2020-03-09 23:58:58" }
)

foreach ($t in $tables) {
    $id = $t.Id
    $sheetName = "!!" + $id
    $ws = $wb.Worksheets.Item($sheetName)

    $newTableLine = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='$id' name='$id' date='$($t.Date)' objTablesVersion='0.0.8'"

    # Sheets ship protected (no password) -- unprotect to edit, then restore.
    $ws.Unprotect()

    if ($id -eq "Compartment") {
        # This first sheet also carries the document-level metadata row (A1),
        # with the table-level row pushed down to A2.
        $ws.Range("A1").Value = "!!!ObjTables schema='SBtab' objTablesVersion='0.0.8' date='$docDate'"
        $ws.Range("A2").Value = $newTableLine
    } else {
        $ws.Range("A1").Value = $newTableLine
    }

    $ws.Protect()
}
